$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Settings sheet: add canvas_item_type setting row (A13/B13)
# ------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("A13").Value = "canvas_item_type"
$settings.Range("B13").Value = "assignment"

# ------------------------------------------------------------------
# Discussions sheet: insert canvas_item_type column (new col L),
# shifting error_message/created_at/updated_at one column right.
# ------------------------------------------------------------------
$disc = $wb.Worksheets.Item("Discussions")

# Insert a new column before the current L (canvas_assignment_id stays
# at K, error_message/created_at/updated_at shift from L/M/N -> M/N/O).
$disc.Columns.Item(12).Insert()

# Header + width for the freshly inserted column L.
$disc.Range("L1").Value = "canvas_item_type"
$disc.Columns.Item(12).ColumnWidth = 19.166666666666668   # stored width -> 20

# The shifted error_message column (now M) should end up width 17.
$disc.Columns.Item(13).ColumnWidth = 16.166666666666668   # stored width -> 17

# Re-create the status (E) validation without the custom error message,
# matching the simplified attribute set in the target file.
$statusRng = $disc.Range("E2:E1001")
$statusRng.Validation.Delete()
$statusRng.Validation.Add(3, 1, 1, '"uploaded,transcribing,mapping,review,approved,sent,error"')
$statusRng.Validation.IgnoreBlank = $false
$statusRng.Validation.InCellDropdown = $true
$statusRng.Validation.ShowInput = $false
$statusRng.Validation.ShowError = $false

# New canvas_item_type (L) validation: assignment vs discussion, blanks allowed.
$itemTypeRng = $disc.Range("L2:L1001")
$itemTypeRng.Validation.Add(3, 1, 1, '"assignment,discussion"')
$itemTypeRng.Validation.IgnoreBlank = $true
$itemTypeRng.Validation.InCellDropdown = $true
$itemTypeRng.Validation.ShowInput = $false
$itemTypeRng.Validation.ShowError = $false
